# Regen save_data: recompute column G (K) values for shaw_bryan.xlsx
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
#                   calc and write s_vals"
# The column G ("K") values were recalculated with a new method and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 3
    11 = 2
    12 = 1
    14 = 1
    16 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
